$p = $ppt.ActivePresentation

# --- Slide 35: "first set for" -> "first check for" ---
$s35 = $p.Slides.Item(35)
$shape35 = $s35.Shapes.Item("Rectangle 4")
$shape35.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "first check for"

# --- Slide 63: merge the first three runs of the second paragraph into one run ---
$s63 = $p.Slides.Item(63)
$shape63 = $s63.Shapes.Item("Content Placeholder 2")
$tr63 = $shape63.TextFrame.TextRange
$full63 = $tr63.Text
$startIdx = $full63.IndexOf("In order to ease")
$endIdx = $full63.IndexOf("in a ") + 5
$sub63 = $tr63.Characters($startIdx + 1, $endIdx - $startIdx)
$sub63.Text = "In order to ease the transition to error recovery in the next version of the parser, most parsing methods will wrap the basic parsing logic in a "

# --- Slide 65: split "parseVariable()" and "parseNamedValue()" runs so the
#     identifier and the trailing "()" are separate runs (identifier keeps the
#     Consolas formatting, matching the new run layout in the diff) ---
$s65 = $p.Slides.Item(65)
$shape65 = $s65.Shapes.Item("Content Placeholder 2")
$tr65 = $shape65.TextFrame.TextRange

$word1 = "parseVariable"
$full65 = $tr65.Text
$idx1 = $full65.IndexOf($word1 + "()")
$sub65a = $tr65.Characters($idx1 + 1, $word1.Length)
$sub65a.Text = $word1

$word2 = "parseNamedValue"
$full65 = $tr65.Text
$idx2 = $full65.IndexOf($word2 + "()")
$sub65b = $tr65.Characters($idx2 + 1, $word2.Length)
$sub65b.Text = $word2
